$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Source template rows: 259 (navigate) and 260 (Close Driver) already carry
# the correct per-cell styles (s="3" on D for navigate-Pass, s="2" across Close Driver).
# Copy that 2-row block repeatedly to extend the log with the new entries, then
# patch just the timestamp (column E) for every new row - every other logged value
# (action/url/status/description) is identical for each new pair.

$ws.Range("A259:G260").Copy($ws.Range("A261"))
$ws.Range("E261").Value = "04/29/2025 14:34:09"
$ws.Range("E262").Value = "04/29/2025 14:34:36"

$ws.Range("A259:G260").Copy($ws.Range("A263"))
$ws.Range("E263").Value = "04/29/2025 14:36:44"
$ws.Range("E264").Value = "04/29/2025 14:37:11"

$ws.Range("A259:G260").Copy($ws.Range("A265"))
$ws.Range("E265").Value = "04/29/2025 15:10:18"
$ws.Range("E266").Value = "04/29/2025 15:10:45"

$ws.Range("A259:G260").Copy($ws.Range("A267"))
$ws.Range("E267").Value = "04/30/2025 10:44:31"
$ws.Range("E268").Value = "04/30/2025 10:44:57"

$ws.Range("A259:G260").Copy($ws.Range("A269"))
$ws.Range("E269").Value = "04/30/2025 10:46:30"
$ws.Range("E270").Value = "04/30/2025 10:46:58"

$ws.Range("A259:G260").Copy($ws.Range("A271"))
$ws.Range("E271").Value = "04/30/2025 10:49:35"
$ws.Range("E272").Value = "04/30/2025 10:49:48"

$ws.Range("A259:G260").Copy($ws.Range("A273"))
$ws.Range("E273").Value = "04/30/2025 10:58:19"
$ws.Range("E274").Value = "04/30/2025 10:58:47"

$ws.Range("A259:G260").Copy($ws.Range("A275"))
$ws.Range("E275").Value = "04/30/2025 11:01:29"
$ws.Range("E276").Value = "04/30/2025 11:01:59"

$ws.Range("A259:G260").Copy($ws.Range("A277"))
$ws.Range("E277").Value = "04/30/2025 11:03:38"
$ws.Range("E278").Value = "04/30/2025 11:04:03"

$ws.Range("A259:G260").Copy($ws.Range("A279"))
$ws.Range("E279").Value = "04/30/2025 11:04:44"
$ws.Range("E280").Value = "04/30/2025 11:05:19"

$ws.Range("A259:G260").Copy($ws.Range("A281"))
$ws.Range("E281").Value = "04/30/2025 11:17:26"
$ws.Range("E282").Value = "04/30/2025 11:18:02"

$ws.Range("A259:G260").Copy($ws.Range("A283"))
$ws.Range("E283").Value = "04/30/2025 11:23:35"
$ws.Range("E284").Value = "04/30/2025 11:24:01"

$ws.Range("A259:G260").Copy($ws.Range("A285"))
$ws.Range("E285").Value = "04/30/2025 11:25:09"
$ws.Range("E286").Value = "04/30/2025 11:25:40"

